# EIA Table 2.3.B (Petroleum Coke: Consumption for Useful Thermal Output)
# Monthly refresh: extend report from "...2006-October 2016" through November 2016.
#
# This:
#   1. Updates the two report-wide header strings that mention the month/year.
#   2. Inserts a new "November" monthly row into the "Year 2016" section
#      (rows 43-52 already hold Jan-Oct 2016; the new row becomes row 53,
#      pushing "Annual Totals" / "Year to Date" / "Rolling 12 Months" blocks
#      down by one row).
#   3. Refreshes the "Year to Date" annual totals (now through November) for
#      2014/2015/2016, and the "Rolling 12 Months" totals for 2015/2016.
#   4. Renames the "Rolling 12 Months Ending in October" label to "...November".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- 1. Report title / subtitle -------------------------------------------
$ws.Cells.Item(2, 1).Value = "by Sector, 2006-November 2016 (Thousand Tons)"

# --- 2. Insert the new "November" row under the Year 2016 monthly block ---
# Row 52 is "October" (the last monthly row of the Year 2016 section);
# inserting above row 53 ("Annual Totals" / "Year to Date" header) makes
# room for November while pushing everything below down by one row.
$ws.Rows.Item(53).Insert()

# Copy the formatting (number formats, borders, fill, etc.) from the row
# above (October) onto the newly inserted blank row so it matches the rest
# of the monthly data rows in this section.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial($xlPasteFormats)

$ws.Cells.Item(53, 1).Value = "November"
$ws.Cells.Item(53, 2).Value = 77
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 9
$ws.Cells.Item(53, 5).Value = 0.16
$ws.Cells.Item(53, 6).Value = 68

# --- 3. Refresh "Year to Date" annual totals (rows shifted down by 1) -----
# Row 54 = "Year to Date" header (was row 53)
# Row 55 = 2014, Row 56 = 2015, Row 57 = 2016
$ws.Cells.Item(55, 2).Value = 1168
$ws.Cells.Item(55, 3).Value = 3
$ws.Cells.Item(55, 4).Value = 80
$ws.Cells.Item(55, 5).Value = 14
$ws.Cells.Item(55, 6).Value = 1071

$ws.Cells.Item(56, 2).Value = 1058
$ws.Cells.Item(56, 3).Value = 9
$ws.Cells.Item(56, 4).Value = 99
$ws.Cells.Item(56, 5).Value = 15
$ws.Cells.Item(56, 6).Value = 935

$ws.Cells.Item(57, 2).Value = 923
$ws.Cells.Item(57, 3).Value = 2
$ws.Cells.Item(57, 4).Value = 95
$ws.Cells.Item(57, 5).Value = 7
$ws.Cells.Item(57, 6).Value = 818

# --- 4. Refresh "Rolling 12 Months" label + totals (rows shifted down) ----
# Row 58 = "Rolling 12 Months Ending in October" header (was row 57)
$ws.Cells.Item(58, 1).Value = "Rolling 12 Months Ending in November"

# Row 59 = 2015, Row 60 = 2016
$ws.Cells.Item(59, 2).Value = 1172
$ws.Cells.Item(59, 3).Value = 9
$ws.Cells.Item(59, 4).Value = 108
$ws.Cells.Item(59, 5).Value = 17
$ws.Cells.Item(59, 6).Value = 1038

$ws.Cells.Item(60, 2).Value = 1009
$ws.Cells.Item(60, 3).Value = 2
$ws.Cells.Item(60, 4).Value = 106
$ws.Cells.Item(60, 5).Value = 9
$ws.Cells.Item(60, 6).Value = 893
